$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 193, shifting the existing 193:201 down to 195:203
$ws.Rows("193:194").Insert()

# Row 193 - new "Primera" quality record
$ws.Cells.Item(193, 1).Value() = 6
$ws.Cells.Item(193, 2).Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(193, 3).Value() = "Metropolitana"
$ws.Cells.Item(193, 4).Value() = 45075
$ws.Cells.Item(193, 5).Value() = 13
$ws.Cells.Item(193, 6).Value() = "Fruta"
$ws.Cells.Item(193, 7).Value() = 100104
$ws.Cells.Item(193, 8).Value() = "Frutos de pepita"
$ws.Cells.Item(193, 9).Value() = 100104003
$ws.Cells.Item(193, 10).Value() = "Membrillo"
$ws.Cells.Item(193, 11).Value() = "Champion"
$ws.Cells.Item(193, 12).Value() = "Primera"
$ws.Cells.Item(193, 13).Value() = 8
$ws.Cells.Item(193, 14).Value() = 190000
$ws.Cells.Item(193, 15).Value() = 190000
$ws.Cells.Item(193, 16).Value() = 190000
$ws.Cells.Item(193, 17).Value() = "$/bins (450 kilos)"
$ws.Cells.Item(193, 18).Value() = "Región de O'Higgins"
$ws.Cells.Item(193, 19).Value() = 422
$ws.Cells.Item(193, 20).Value() = 450

# Row 194 - new "Segunda" quality record
$ws.Cells.Item(194, 1).Value() = 6
$ws.Cells.Item(194, 2).Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(194, 3).Value() = "Metropolitana"
$ws.Cells.Item(194, 4).Value() = 45075
$ws.Cells.Item(194, 5).Value() = 13
$ws.Cells.Item(194, 6).Value() = "Fruta"
$ws.Cells.Item(194, 7).Value() = 100104
$ws.Cells.Item(194, 8).Value() = "Frutos de pepita"
$ws.Cells.Item(194, 9).Value() = 100104003
$ws.Cells.Item(194, 10).Value() = "Membrillo"
$ws.Cells.Item(194, 11).Value() = "Champion"
$ws.Cells.Item(194, 12).Value() = "Segunda"
$ws.Cells.Item(194, 13).Value() = 10
$ws.Cells.Item(194, 14).Value() = 160000
$ws.Cells.Item(194, 15).Value() = 160000
$ws.Cells.Item(194, 16).Value() = 160000
$ws.Cells.Item(194, 17).Value() = "$/bins (450 kilos)"
$ws.Cells.Item(194, 18).Value() = "Región de O'Higgins"
$ws.Cells.Item(194, 19).Value() = 356
$ws.Cells.Item(194, 20).Value() = 450

Write-Output "edit complete"
